$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear H2 (previously TRUE, now removed entirely)
$ws.Range("H2").ClearContents()

# Populate header row and data rows with updated / new values
$ws.Range("A1").Value = 'BrowserFileDialog'
$ws.Range("B1").Value = 'JScript1 variable'
$ws.Range("C1").Value = 'JScript1 message'
$ws.Range("D1").Value = 'Compression'
$ws.Range("E1").Value = 'Cache'
$ws.Range("F1").Value = 'Profiler'
$ws.Range("G1").Value = 'Python server'
$ws.Range("H1").Value = 'File Sent'
$ws.Range("I1").Value = 'GChrome'
$ws.Range("J1").Value = 'FireFox'
$ws.Range("K1").Value = 'Build time'
$ws.Range("L1").Value = 'Comment'
$ws.Range("A2").Value = 'Debug.Log(fileurl)'
$ws.Range("B2").Value = $true
$ws.Range("C2").Value = '"1"; URL; "5"'
$ws.Range("D2").Value = $false
$ws.Range("E2").Value = $true
$ws.Range("F2").Value = $true
$ws.Range("G2").Value = $true
$ws.Range("I2").Value = $true
$ws.Range("J2").Value = $true
$ws.Range("A3").Value = 'Debug.Log(fileurl)'
$ws.Range("B3").Value = $true
$ws.Range("C3").Value = '"1"; URL; "5"'
$ws.Range("D3").Value = 'GZIP'
$ws.Range("E3").Value = $true
$ws.Range("F3").Value = $false
$ws.Range("G3").Value = $true
$ws.Range("H3").Value = 'Test.txt'
$ws.Range("I3").Value = $true
$ws.Range("J3").Value = $true
$ws.Range("K3").Value = 16
$ws.Range("A4").Value = 'Debug.Log(fileurl)'
$ws.Range("B4").Value = $true
$ws.Range("C4").Value = 'URL'
$ws.Range("D4").Value = 'GZIP'
$ws.Range("E4").Value = $true
$ws.Range("F4").Value = $false
$ws.Range("G4").Value = $true
$ws.Range("H4").Value = 'Test.txt'
$ws.Range("I4").Value = $true
$ws.Range("J4").Value = $true
$ws.Range("K4").Value = 5.5
$ws.Range("A5").Value = 'Debug.Log(fileurl)'
$ws.Range("B5").Value = $true
$ws.Range("C5").Value = 'URL'
$ws.Range("D5").Value = 'GZIP'
$ws.Range("E5").Value = $true
$ws.Range("F5").Value = $false
$ws.Range("G5").Value = $true
$ws.Range("H5").Value = 'The Mom Test.docx'
$ws.Range("I5").Value = $true
$ws.Range("J5").Value = $true
$ws.Range("A6").Value = 'Debug.Log(fileurl)'
$ws.Range("B6").Value = $true
$ws.Range("C6").Value = 'URL'
$ws.Range("D6").Value = 'GZIP'
$ws.Range("E6").Value = $true
$ws.Range("F6").Value = $false
$ws.Range("G6").Value = $true
$ws.Range("H6").Value = 'Carp DS (csv)'
$ws.Range("I6").Value = $true
$ws.Range("J6").Value = $true
$ws.Range("A7").Value = 'Debug.Log(UrlTextField.text)'
$ws.Range("B7").Value = $true
$ws.Range("C7").Value = 'URL'
$ws.Range("D7").Value = 'GZIP'
$ws.Range("E7").Value = $true
$ws.Range("F7").Value = $false
$ws.Range("G7").Value = $true
$ws.Range("H7").Value = 'Test.txt'
$ws.Range("I7").Value = $false
$ws.Range("J7").Value = $false
$ws.Range("K7").Value = 6
$ws.Range("A8").Value = 'Full coroutine; UrlTextField.text disabled'
$ws.Range("B8").Value = $true
$ws.Range("C8").Value = 'URL'
$ws.Range("D8").Value = 'GZIP'
$ws.Range("E8").Value = $true
$ws.Range("F8").Value = $false
$ws.Range("G8").Value = $true
$ws.Range("H8").Value = 'Test.txt'
$ws.Range("I8").Value = $true
$ws.Range("J8").Value = $true
$ws.Range("K8").Value = 8.25
$ws.Range("A9").Value = 'Full coroutine; UrlTextField.text disabled'
$ws.Range("B9").Value = $true
$ws.Range("C9").Value = 'URL'
$ws.Range("D9").Value = 'GZIP'
$ws.Range("E9").Value = $true
$ws.Range("F9").Value = $false
$ws.Range("G9").Value = $true
$ws.Range("H9").Value = 'Test.csv'
$ws.Range("I9").Value = $true
$ws.Range("J9").Value = $true
$ws.Range("A10").Value = 'Full coroutine; UrlTextField.text disabled'
$ws.Range("B10").Value = $true
$ws.Range("C10").Value = 'URL'
$ws.Range("D10").Value = 'GZIP'
$ws.Range("E10").Value = $true
$ws.Range("F10").Value = $false
$ws.Range("G10").Value = $true
$ws.Range("H10").Value = 'Carp DS (csv)'
$ws.Range("I10").Value = $true
$ws.Range("J10").Value = $true
$ws.Range("L10").Value = 'Ran out of memory in attempting to print'
$ws.Range("A11").Value = 'Coroutine w/out printing file contents; UrlTextField.text disabled'
$ws.Range("B11").Value = $true
$ws.Range("C11").Value = 'URL'
$ws.Range("D11").Value = 'GZIP'
$ws.Range("E11").Value = $true
$ws.Range("F11").Value = $false
$ws.Range("G11").Value = $true
$ws.Range("H11").Value = 'Carp DS (csv)'
$ws.Range("I11").Value = $true
$ws.Range("J11").Value = $true
$ws.Range("K11").Value = 6.75
$ws.Range("A12").Value = 'Coroutine; prints type (string) and length (15M)'
$ws.Range("B12").Value = $true
$ws.Range("C12").Value = 'URL'
$ws.Range("D12").Value = 'GZIP'
$ws.Range("E12").Value = $true
$ws.Range("F12").Value = $false
$ws.Range("G12").Value = $true
$ws.Range("H12").Value = 'Carp DS (csv)'
$ws.Range("I12").Value = $true
$ws.Range("J12").Value = $true
$ws.Range("K12").Value = 6.5

# Center-align the GZIP column (D3:D12) per new style
$ws.Range("D3:D12").HorizontalAlignment = -4108

# Adjust column widths: A wider, B:J uniform width matching existing data columns
$ws.Columns.Item(1).ColumnWidth = 29.75
$ws.Columns.Item(10).ColumnWidth = 15.67

# Update selected cell shown in the saved view
$ws.Range("C14").Select()
